$d = $word.ActiveDocument

# 1. Update the date/time text in the Date-styled paragraph.
$d.Content.Find.Execute("June   1, 2021 (01:53:26 AM)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "June   1, 2021 (05:46:13 PM)", 2)

# 2. Extend "your operation." with the new MacOS-shortcut explanation.
$full = $d.Content
$full.Find.Execute("your operation.") | Out-Null
$full.Collapse(0)                       # collapse to just after "your operation."

$full.InsertAfter(" (")
$full.Collapse(0)

# "This shortcut is for Windows" needs italic formatting.
$italicStart = $full.End
$full.InsertAfter("This shortcut is for Windows")
$italicEnd = $full.End
$italicRange = $d.Range($italicStart, $italicEnd)
$italicRange.Font.Italic = 1

$full.Collapse(0)
$full.InsertAfter("; for MacOS, to undo your operation, use")
$full.Collapse(0)
$full.InsertAfter(" ")
$full.Collapse(0)
$full.InsertAfter("CMD")
$full.Collapse(0)
$full.InsertAfter(" ")
$full.Collapse(0)
$full.InsertAfter("+")
$full.Collapse(0)
$full.InsertAfter(" ")
$full.Collapse(0)
$full.InsertAfter("z")
$full.Collapse(0)
$full.InsertAfter(" ")
$full.Collapse(0)
$full.InsertAfter("instead).")
